$d = $word.ActiveDocument

# --- 1. Split the first paragraph's text into the original run (with two
#        trailing spaces) plus three red "(This is a change - Version for
#        main branch)" runs. ---
$p1 = $d.Paragraphs(1)

$pos = $p1.Range.End - 1
$seg = $d.Range($pos, $pos)
$seg.InsertAfter("  ")

$pos = $p1.Range.End - 1
$seg = $d.Range($pos, $pos)
$seg.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$seg.Font.Color = 255

$pos = $p1.Range.End - 1
$seg = $d.Range($pos, $pos)
$seg.InsertAfter("rsion for main branch")
$seg.Font.Color = 255

$pos = $p1.Range.End - 1
$seg = $d.Range($pos, $pos)
$seg.InsertAfter(")")
$seg.Font.Color = 255

# --- 2. Append a new, empty, shaded paragraph after the final paragraph
#        (just before the section break). ---
$endPos = $d.Content.End
$tail = $d.Range($endPos, $endPos)
$tail.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')
